# Updating odds values in the Betfair Back/Lay workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.75
$ws.Range("P2").Value = 2.9
$ws.Range("Q2").Value = 1.47
$ws.Range("R2").Value = 1.77
$ws.Range("S2").Value = 2.18
$ws.Range("T2").Value = 1.54
$ws.Range("U2").Value = 2.7
$ws.Range("AI2").Value = 44
$ws.Range("N3").Value = 4.7
$ws.Range("Q3").Value = 1.67
$ws.Range("S3").Value = 2.68
$ws.Range("T3").Value = 1.65
$ws.Range("X3").Value = 24
$ws.Range("Y3").Value = 20
$ws.Range("AA3").Value = 100
$ws.Range("AI3").Value = 60
$ws.Range("AL3").Value = 32
$ws.Range("AM3").Value = 75
$ws.Range("AO3").Value = 40
$ws.Range("Y4").Value = 1000
$ws.Range("F5").Value = 1.45
$ws.Range("G5").Value = 1.85
$ws.Range("H5").Value = 5
$ws.Range("J5").Value = 3.35
$ws.Range("K5").Value = 9
$ws.Range("L5").Value = 1.31
$ws.Range("N5").Value = 2.7
$ws.Range("P5").Value = 1.7
$ws.Range("Q5").Value = 1.8
$ws.Range("R5").Value = 1.22
$ws.Range("S5").Value = 2.62
$ws.Range("W5").Value = 2.16
$ws.Range("F6").Value = 13.5
$ws.Range("G6").Value = 19.5
$ws.Range("H6").Value = 1.2
$ws.Range("I6").Value = 1.31
$ws.Range("J6").Value = 6.2
$ws.Range("K6").Value = 8.6
$ws.Range("F7").Value = 2
$ws.Range("I7").Value = 5.1
$ws.Range("J7").Value = 2.6
$ws.Range("N7").Value = 1.62
$ws.Range("P7").Value = 1.62
$ws.Range("R7").Value = 1.2
$ws.Range("V7").Value = 1.24
$ws.Range("Q8").Value = 2.16
$ws.Range("R8").Value = 1.31
$ws.Range("S8").Value = 4
$ws.Range("T9").Value = 1.77
$ws.Range("AG9").Value = 9.4
$ws.Range("AJ9").Value = 14.5
$ws.Range("Q11").Value = 2.02
$ws.Range("AB11").Value = 10
$ws.Range("F12").Value = 8.6
$ws.Range("G12").Value = 8.800000000000001
$ws.Range("N12").Value = 5.1
$ws.Range("T12").Value = 1.95
$ws.Range("W12").Value = 1.12
$ws.Range("Z12").Value = 8.6
$ws.Range("AN12").Value = 130
